$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 204; this shifts existing rows 204-246 down to 205-247
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new weekly record
$ws.Range("A204").Value = 4
$ws.Range("B204").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value = "Los Lagos"
$ws.Range("D204").Value = 44782
$ws.Range("E204").Value = 10
$ws.Range("F204").Value = 100112039
$ws.Range("G204").Value = "Ciboulette"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 240
$ws.Range("K204").Value = 4000
$ws.Range("L204").Value = 4000
$ws.Range("M204").Value = 4000
$ws.Range("N204").Value = "$/docena de atados"
$ws.Range("O204").Value = "Región Metropolitana"
$ws.Range("P204").Value = 1333
$ws.Range("Q204").Value = 3
$ws.Range("R204").Value = "Hortaliza"
